$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Range("B5").Value = 86
$ws.Range("B7").Value = 89
$ws.Range("B8").Value = 91
$ws.Range("B10").Value = 94
$ws.Range("B11").Value = 96
$ws.Range("B12").Value = 98
$ws.Range("B13").Value = 100
$ws.Range("B14").Value = 102
$ws.Range("B15").Value = 104
$ws.Range("B21").Value = 119
$ws.Range("B22").Value = 122
$ws.Range("B23").Value = 126

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Range("B2").Value = 73
$ws.Range("B12").Value = 91
$ws.Range("B13").Value = 93
$ws.Range("B14").Value = 95

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Range("B3").Value = 68
$ws.Range("B5").Value = 71
$ws.Range("B17").Value = 94
$ws.Range("B19").Value = 99
$ws.Range("B21").Value = 105
$ws.Range("B22").Value = 108
$ws.Range("B23").Value = 112
$ws.Range("B25").Value = 123
$ws.Range("B26").Value = 130

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Range("B2").Value = 60
$ws.Range("B3").Value = 62
$ws.Range("B4").Value = 63
$ws.Range("B5").Value = 65
$ws.Range("B6").Value = 66
$ws.Range("B7").Value = 68
$ws.Range("B8").Value = 70
$ws.Range("B9").Value = 72
$ws.Range("B10").Value = 73
$ws.Range("B11").Value = 75
$ws.Range("B12").Value = 77
$ws.Range("B13").Value = 79
$ws.Range("B14").Value = 81
$ws.Range("B15").Value = 83
$ws.Range("B16").Value = 86
$ws.Range("B17").Value = 88
$ws.Range("B18").Value = 90
$ws.Range("B19").Value = 93
$ws.Range("B20").Value = 96
$ws.Range("B21").Value = 99
$ws.Range("B22").Value = 102
$ws.Range("B23").Value = 106
$ws.Range("B24").Value = 111
$ws.Range("B25").Value = 117
$ws.Range("B26").Value = 126

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Range("B2").Value = 54
$ws.Range("B3").Value = 55
$ws.Range("B4").Value = 57
$ws.Range("B5").Value = 58
$ws.Range("B11").Value = 69
$ws.Range("B15").Value = 77
$ws.Range("B25").Value = 115

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Range("B2").Value = 48
$ws.Range("B3").Value = 49
$ws.Range("B4").Value = 51
$ws.Range("B5").Value = 52
$ws.Range("B6").Value = 54
$ws.Range("B7").Value = 55
$ws.Range("B9").Value = 58
$ws.Range("B14").Value = 68
$ws.Range("B21").Value = 87
$ws.Range("B24").Value = 107
$ws.Range("B25").Value = 118
$ws.Range("B26").Value = 124
$ws.Range("B27").Value = 129
